# funcion de lectura de archivos excel
# The "payment_id" column (A) is no longer needed for the read workflow,
# so it is dropped and the remaining columns (sale_id, payment_date,
# amount, payment_method) shift one slot to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A ("payment_id") entirely; B:E shifts left into A:D,
# carrying its values, shared-string types and cell styles with it.
$ws.Range("A1").EntireColumn.Delete()

# Leave the selection where the author's cursor ended up.
$ws.Range("C10").Select()
